$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 (b.md) status updated to "Ready for handoff" ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-17 04:33:44"

# --- zh-cn sheet: row 3 (b.md) handoff info updated ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C3").Value = "Ready for handoff"

# "False" looks boolean to Excel's auto-detection; force text with a leading
# apostrophe, then restore the plain (non quote-prefixed) style copied from
# the cell above so no stray formatting is introduced.
$plainStyleZh = $ws2.Range("F2").Style
$ws2.Range("F3").Value = "'False"
$ws2.Range("F3").Style = $plainStyleZh

$ws2.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-17 04:33:40"
$ws2.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ffb1665bb9d0c49f3b4fe5930ab1892091a09130/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fe8e98219f8c89ae630a5ad01729b86a41b1ad3/e2e/b.md."

# --- de-de sheet: row 3 (b.md) handoff info updated ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C3").Value = "Ready for handoff"

$plainStyleDe = $ws3.Range("F2").Style
$ws3.Range("F3").Value = "'False"
$ws3.Range("F3").Style = $plainStyleDe

$ws3.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-17 04:33:44"
$ws3.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ffb1665bb9d0c49f3b4fe5930ab1892091a09130/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fe8e98219f8c89ae630a5ad01729b86a41b1ad3/e2e/b.md."

# --- widen column P (Error Detail) on both locale sheets to fit the new message ---
# (39.17 is the Excel "characters" ColumnWidth that round-trips to an OOXML
# <col .../> width of exactly 40, matching columns G/J which already use it.)
$ws2.Columns.Item(16).ColumnWidth = 39.17
$ws3.Columns.Item(16).ColumnWidth = 39.17

Write-Host "Applied handoff report updates"
